$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Insert a new row at 76 (shifts rows 76-104 down to 77-105) to make room
# for the new "pidSV(<int>)" row that keeps the original pidSV description.
$ws.Rows.Item(76).Insert()

# Row 75 becomes the new "adjustSV" command (replacing the old pidSV row).
$ws.Range("B75").Value = "adjustSV(<int>)"
$ws.Range("C75").Value = "increases or decreases the current target SV value by <int>"

# New row 76 keeps the original pidSV text/description, now with <int> instead of <float>.
$ws.Range("B76").Value = "pidSV(<int>)"
$ws.Range("C76").Value = "sets the PID target set value SV"
$ws.Rows.Item(76).RowHeight = 13.8
